$wb = $excel.ActiveWorkbook

# --- Stocks sheet: add a new Apple (USD) row before the PF line, and keep a
#     trailing blank (dated) row available under it -----------------------
$stocks = $wb.Worksheets.Item("Stocks")

$stocks.Range("A4").Value = 44062
$stocks.Range("B4").Value = "Apple"
$stocks.Range("C4").Value = "AAPL"
$stocks.Range("D4").Value = 1
$stocks.Range("E4").Value = 440
$stocks.Range("F4").Value = 372.52
$stocks.Range("G4").Value = "USD"
$stocks.Range("H4").Value = 0.5
$stocks.Range("I4").Value = "#C0C0C0"
$stocks.Range("J4").Value = $false
$stocks.Range("K4").Value = "Tech"

# Grow the table with a fresh blank (date-formatted) row, mirroring the
# empty placeholder row already sitting above it.
$stocks.Range("A5").Copy($stocks.Range("A6"))

$stocks.Range("D19").Select()

# --- Interests sheet: this stray duplicate (no more constant scraping of
#     this quote) moves to Stocks above, so drop it here ------------------
$interests = $wb.Worksheets.Item("Interests")

$interests.Range("A9:K9").ClearContents()

$interests.Range("A6:K8").Select()

# Leave "Stocks" as the active/visible tab, matching the saved workbook state.
$stocks.Activate()
